# Change end game condition and update strings to match.
#
# The game used to end when too many cases went unresolved. Now the game
# ends when citizen happiness falls too low. This updates the "game over"
# strings accordingly and adds a new "Citizen Happiness" localisation key,
# mirroring the row layout already used for the other string entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "BASIC_TEXT_GAMEOVER" (row 69) - the game-over headline text.
$ws.Range("B69").Value = "Game Over*2n*Citizens not happy with your performance"

# "BASIC_TEXT_GAMEOVER_BODY" (row 70) - the game-over body text.
$ws.Range("B70").Value = "You Survived {0} Turns*2n*Before citizen happiness fell too low, Citizens no longer feel safe under your control"

# New row appended after "BASIC_TEXT_OFFICERS_REQUIRED" (row 73) for the new
# citizen-happiness localisation key.
$ws.Range("A74").Value = "BASIC_TEXT_CITIZEN_HAPPINESS"
$ws.Range("B74").Value = "Citizen Happines"
$ws.Range("C74").Value = "XXXX"
$ws.Range("D74").Value = "XXXX"
$ws.Range("E74").Value = "XXXX"

# Update the saved view state: scroll so row 37 is at the top, and the
# active/selected cell is B60.
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 1
$ws.Range("B60").Select()
